$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $origStyle = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $origStyle
}

Set-TextValue $ws "D2" '67.329.47'
$ws.Range("E2").Value = '  -0.19%  '
Set-TextValue $ws "D3" '3.501.73'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws "D5" '598.90'
$ws.Range("E5").Value = '  +0.28%  '
Set-TextValue $ws "D6" '175.44'
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -0.82%  '
$ws.Range("E9").Value = '  -2.88%  '
Set-TextValue $ws "D10" '7.14'
$ws.Range("E10").Value = '  -2.79%  '
Set-TextValue $ws "D11" '0.431'
$ws.Range("E11").Value = '  -1.88%  '
Set-TextValue $ws "D12" '4.111.65'
$ws.Range("E12").Value = '  -0.69%  '
Set-TextValue $ws "D13" '31.13'
$ws.Range("E13").Value = '  +8.11%  '
$ws.Range("E14").Value = '  -0.01%  '
Set-TextValue $ws "D15" '67.353.85'
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("E16").Value = '  -2.27%  '
Set-TextValue $ws "D17" '3.504.90'
$ws.Range("E17").Value = '  -0.56%  '
Set-TextValue $ws "D18" '6.30'
$ws.Range("E18").Value = '  -1.47%  '
Set-TextValue $ws "D19" '14.60'
$ws.Range("E19").Value = '  +2.31%  '
Set-TextValue $ws "D20" '391.98'
$ws.Range("E20").Value = '  -1.48%  '
Set-TextValue $ws "D21" '7.97'
$ws.Range("E21").Value = '  -0.90%  '
Set-TextValue $ws "D22" '73.51'
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws "D23" '0.999'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws "D24" '0.540'
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  +0.11%  '
Set-TextValue $ws "D26" '0.0000122'
$ws.Range("E26").Value = '  -1.81%  '
Set-TextValue $ws "D27" '10.28'
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("E28").Value = '  -1.13%  '
Set-TextValue $ws "D29" '0.994'
$ws.Range("E29").Value = '  -0.64%  '
Set-TextValue $ws "D30" '6.12'
$ws.Range("E30").Value = '  -3.76%  '
Set-TextValue $ws "D31" '1.42'
$ws.Range("E31").Value = '  -4.07%  '
$ws.Range("E32").Value = '  -0.76%  '
Set-TextValue $ws "D33" '23.65'
$ws.Range("E33").Value = '  -2.24%  '
Set-TextValue $ws "D34" '7.38'
$ws.Range("E34").Value = '  -1.17%  '
Set-TextValue $ws "D35" '1.64'
$ws.Range("E35").Value = '  +0.39%  '
Set-TextValue $ws "D36" '163.27'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws "D37" '0.878'
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws "D38" '1.94'
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("E39").Value = '  +1.61%  '
Set-TextValue $ws "D40" '4.66'
$ws.Range("E40").Value = '  -2.46%  '
Set-TextValue $ws "D41" '26.50'
$ws.Range("E41").Value = '  -1.26%  '
Set-TextValue $ws "D42" '27.11'
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws "D43" '2.812.55'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws "D44" '0.0730'
$ws.Range("E44").Value = '  -2.98%  '
Set-TextValue $ws "D45" '42.56'
$ws.Range("E45").Value = '  -1.07%  '
Set-TextValue $ws "D46" '2.53'
$ws.Range("E46").Value = '  -4.11%  '
$ws.Range("E47").Value = '  -4.51%  '
Set-TextValue $ws "D48" '337.09'
$ws.Range("E48").Value = '  -1.91%  '
$ws.Range("E49").Value = '  -2.88%  '
Set-TextValue $ws "D50" '33.65'
$ws.Range("E50").Value = '  -0.38%  '
Set-TextValue $ws "D51" '0.846'
$ws.Range("E51").Value = '  -2.00%  '

Write-Host "Applied cryptos list update"
